$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between row 70 and row 71 ---
# New row 70 (was row 71's match): Famalicao vs Vitoria Guimaraes
$ws.Cells.Item(70, 6).Value = "Famalicao"
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = "Vitoria Guimaraes"
$ws.Cells.Item(70, 9).Value = 3
$ws.Cells.Item(70, 10).Value = 2.59
$ws.Cells.Item(70, 11).Value = "02/10/2023 07:12"
$ws.Cells.Item(70, 12).Value = 2.94
$ws.Cells.Item(70, 13).Value = "08/10/2023 16:27"
$ws.Cells.Item(70, 14).Value = 3.26
$ws.Cells.Item(70, 15).Value = "02/10/2023 07:12"
$ws.Cells.Item(70, 16).Value = 3.07
$ws.Cells.Item(70, 17).Value = "08/10/2023 16:27"
$ws.Cells.Item(70, 18).Value = 2.89
$ws.Cells.Item(70, 19).Value = "02/10/2023 07:12"
$ws.Cells.Item(70, 20).Value = 2.77
$ws.Cells.Item(70, 21).Value = "08/10/2023 16:27"
$ws.Cells.Item(70, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/famalicao-vitoria-guimaraes/beAOJg87/"

# New row 71 (was row 70's match): Casa Pia vs Estrela
$ws.Cells.Item(71, 6).Value = "Casa Pia"
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = "Estrela"
$ws.Cells.Item(71, 9).Value = 1
$ws.Cells.Item(71, 10).Value = 1.89
$ws.Cells.Item(71, 11).Value = "02/10/2023 20:42"
$ws.Cells.Item(71, 12).Value = 2.16
$ws.Cells.Item(71, 13).Value = "08/10/2023 16:29"
$ws.Cells.Item(71, 14).Value = 3.53
$ws.Cells.Item(71, 15).Value = "02/10/2023 20:42"
$ws.Cells.Item(71, 16).Value = 3.47
$ws.Cells.Item(71, 17).Value = "08/10/2023 16:29"
$ws.Cells.Item(71, 18).Value = 4.36
$ws.Cells.Item(71, 19).Value = "02/10/2023 20:42"
$ws.Cells.Item(71, 20).Value = 3.61
$ws.Cells.Item(71, 21).Value = "08/10/2023 16:28"
$ws.Cells.Item(71, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/casa-pia-estrela-da-amadora/Cbb6rwo8/"

# --- Swap match data (columns F:V) between row 76 and row 77 ---
# New row 76 (was row 77's match): Benfica vs Casa Pia
$ws.Cells.Item(76, 6).Value = "Benfica"
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = "Casa Pia"
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 1.22
$ws.Cells.Item(76, 11).Value = "11/10/2023 14:42"
$ws.Cells.Item(76, 12).Value = 1.22
$ws.Cells.Item(76, 13).Value = "28/10/2023 18:55"
$ws.Cells.Item(76, 14).Value = 7.32
$ws.Cells.Item(76, 15).Value = "11/10/2023 14:42"
$ws.Cells.Item(76, 16).Value = 6.95
$ws.Cells.Item(76, 17).Value = "28/10/2023 18:58"
$ws.Cells.Item(76, 18).Value = 13.22
$ws.Cells.Item(76, 19).Value = "11/10/2023 14:42"
$ws.Cells.Item(76, 20).Value = 14.17
$ws.Cells.Item(76, 21).Value = "28/10/2023 18:58"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/benfica-casa-pia/GWtkzFhl/"

# New row 77 (was row 76's match): Vitoria Guimaraes vs Chaves
$ws.Cells.Item(77, 6).Value = "Vitoria Guimaraes"
$ws.Cells.Item(77, 7).Value = 5
$ws.Cells.Item(77, 8).Value = "Chaves"
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 1.62
$ws.Cells.Item(77, 11).Value = "11/10/2023 14:42"
$ws.Cells.Item(77, 12).Value = 1.69
$ws.Cells.Item(77, 13).Value = "28/10/2023 18:58"
$ws.Cells.Item(77, 14).Value = 4.26
$ws.Cells.Item(77, 15).Value = "11/10/2023 14:42"
$ws.Cells.Item(77, 16).Value = 4
$ws.Cells.Item(77, 17).Value = "28/10/2023 18:58"
$ws.Cells.Item(77, 18).Value = 5.63
$ws.Cells.Item(77, 19).Value = "11/10/2023 14:42"
$ws.Cells.Item(77, 20).Value = 5.29
$ws.Cells.Item(77, 21).Value = "28/10/2023 18:58"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/vitoria-guimaraes-chaves/8vH9wlat/"

# --- Append new row 92 (Estoril vs Casa Pia) ---
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V92").PasteSpecial(-4122)

$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "portugal"
$ws.Cells.Item(92, 3).Value = "liga-portugal"
$ws.Cells.Item(92, 4).Value = "2023-2024"
$ws.Cells.Item(92, 5).Value = 45240.88541666666
$ws.Cells.Item(92, 6).Value = "Estoril"
$ws.Cells.Item(92, 7).Value = 4
$ws.Cells.Item(92, 8).Value = "Casa Pia"
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 2.49
$ws.Cells.Item(92, 11).Value = "05/11/2023 16:42"
$ws.Cells.Item(92, 12).Value = 2.46
$ws.Cells.Item(92, 13).Value = "10/11/2023 21:07"
$ws.Cells.Item(92, 14).Value = 3.32
$ws.Cells.Item(92, 15).Value = "05/11/2023 16:42"
$ws.Cells.Item(92, 16).Value = 3.31
$ws.Cells.Item(92, 17).Value = "10/11/2023 21:07"
$ws.Cells.Item(92, 18).Value = 3.04
$ws.Cells.Item(92, 19).Value = "05/11/2023 16:42"
$ws.Cells.Item(92, 20).Value = 3.14
$ws.Cells.Item(92, 21).Value = "10/11/2023 21:07"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/estoril-casa-pia/ppfZ24Si/"
